$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("AQ2").Value = 63
$ws.Range("BH2").Value = 255
$ws.Range("AW4").Value = 233
$ws.Range("BH5").Value = 432
$ws.Range("AR6").Value = 1312
$ws.Range("BH6").Value = 342
$ws.Range("BL6").Value = 1061
$ws.Range("BO6").Value = 756
$ws.Range("AQ7").Value = 639
$ws.Range("AR7").Value = 411
$ws.Range("AX7").Value = 220
$ws.Range("BC7").Value = 395
$ws.Range("BH7").Value = 58
$ws.Range("BJ7").Value = 509
$ws.Range("BL7").Value = 437
$ws.Range("BO7").Value = 550
$ws.Range("BT7").Value = 606
$ws.Range("BG9").Value = 1077
$ws.Range("BL9").Value = 343
$ws.Range("BM9").Value = 500
$ws.Range("BO9").Value = 131
$ws.Range("AR11").Value = 633
$ws.Range("AW11").Value = 275
$ws.Range("BG11").Value = 1805
$ws.Range("BJ11").Value = 631
$ws.Range("BL11").Value = 401
$ws.Range("BO11").Value = 709
$ws.Range("BT12").Value = 136
$ws.Range("AR13").Value = 1427
$ws.Range("BL13").Value = 398
$ws.Range("BO13").Value = 1651
$ws.Range("BT13").Value = 1167
$ws.Range("AW14").Value = 2270
$ws.Range("BL14").Value = 2015
$ws.Range("BO14").Value = 2041
$ws.Range("BT14").Value = 1479
$ws.Range("AR15").Value = 6425
$ws.Range("AW15").Value = 4106
$ws.Range("AX15").Value = 6962
$ws.Range("BH15").Value = 4414
$ws.Range("BJ15").Value = 3794
$ws.Range("BL15").Value = 2548
$ws.Range("BM15").Value = 3619
$ws.Range("BO15").Value = 2169
$ws.Range("BT15").Value = 2063
$ws.Range("BU15").Value = 2640
$ws.Range("AQ16").Value = 1351
$ws.Range("AR16").Value = 2288
$ws.Range("AW16").Value = 1248
$ws.Range("AX16").Value = 870
$ws.Range("BH16").Value = 916
$ws.Range("BJ16").Value = 1202
$ws.Range("BL16").Value = 2396
$ws.Range("BM16").Value = 452
$ws.Range("BO16").Value = 1379
$ws.Range("AR17").Value = 6742
$ws.Range("AW17").Value = 3051
$ws.Range("AX17").Value = 3464
$ws.Range("BC17").Value = 1228
$ws.Range("BH17").Value = 687
$ws.Range("BJ17").Value = 2955
$ws.Range("BL17").Value = 2547
$ws.Range("BM17").Value = 2877
$ws.Range("BO17").Value = 2117
$ws.Range("BT17").Value = 1322
$ws.Range("AW20").Value = 7124
$ws.Range("AX20").Value = 7938
$ws.Range("BG20").Value = 3523
$ws.Range("BH20").Value = 3189
$ws.Range("BL20").Value = 8318
$ws.Range("BM20").Value = 4245
$ws.Range("BO20").Value = 7030
$ws.Range("BT20").Value = 6531
$ws.Range("AR21").Value = 224
$ws.Range("AW21").Value = 714
$ws.Range("AX21").Value = 180
$ws.Range("BC21").Value = 139
$ws.Range("BL21").Value = 107
$ws.Range("AR22").Value = 11338
$ws.Range("AW22").Value = 7073
$ws.Range("AX22").Value = 5518
$ws.Range("BG22").Value = 3299
$ws.Range("BL22").Value = 7345
$ws.Range("BM22").Value = 4102
$ws.Range("BO22").Value = 6586
$ws.Range("BT22").Value = 4280
$ws.Range("BU22").Value = 1768
$ws.Range("BC23").Value = 902
$ws.Range("BC24").Value = 209
$ws.Range("BJ24").Value = 346
$ws.Range("AQ25").Value = 5052
$ws.Range("BO25").Value = 3418
$ws.Range("AX26").Value = 1499
$ws.Range("BO31").Value = 105
$ws.Range("BT32").Value = 71
$ws.Range("AR34").Value = 181
$ws.Range("BC34").Value = 258
$ws.Range("BO34").Value = 371
$ws.Range("BT34").Value = 124
$ws.Range("BT35").Value = 140
$ws.Range("BH37").Value = 157
$ws.Range("BC38").Value = 656
$ws.Range("BH38").Value = 17
$ws.Range("BO38").Value = 1393
$ws.Range("BT38").Value = 688
$ws.Range("BU38").Value = 746
$ws.Range("BT39").Value = 559
$ws.Range("AR40").Value = 1048
$ws.Range("BC40").Value = 413
$ws.Range("BG40").Value = 8069
$ws.Range("BH40").Value = 744
$ws.Range("BO40").Value = 349
$ws.Range("BT40").Value = 410
$ws.Range("BH41").Value = 487
$ws.Range("BL41").Value = 1652
$ws.Range("BO41").Value = 1488
$ws.Range("BT41").Value = 1517
$ws.Range("BT42").Value = 59
$ws.Range("AX44").Value = 295
$ws.Range("BC44").Value = 160
$ws.Range("BH44").Value = 118
$ws.Range("BL44").Value = 405
$ws.Range("BM44").Value = 191
$ws.Range("BO44").Value = 610
$ws.Range("AX45").Value = 2787
$ws.Range("BL45").Value = 1594
$ws.Range("BJ48").Value = 297
$ws.Range("BL48").Value = 259
$ws.Range("BM48").Value = 497
$ws.Range("BT48").Value = 412
$ws.Range("BG49").Value = 3111
$ws.Range("BL50").Value = 620
$ws.Range("BM50").Value = 471
$ws.Range("BT50").Value = 440
$ws.Range("BT53").Value = 366
$ws.Range("BH56").Value = 416
$ws.Range("BJ56").Value = 2493
$ws.Range("AW60").Value = 109
$ws.Range("BC60").Value = 122
$ws.Range("BL62").Value = 697
$ws.Range("BO62").Value = 1591
$ws.Range("BT62").Value = 1533
$ws.Range("BL64").Value = 1024
$ws.Range("BO64").Value = 427
$ws.Range("BT64").Value = 484
$ws.Range("BH67").Value = 154
$ws.Range("BJ67").Value = 311
$ws.Range("AW68").Value = 134
$ws.Range("AX68").Value = 438
$ws.Range("BT68").Value = 206
$ws.Range("BT69").Value = 519
$ws.Range("BL70").Value = 160
$ws.Range("BO71").Value = 1051
$ws.Range("BO73").Value = 932
$ws.Range("BO74").Value = 252
$ws.Range("BO75").Value = 414
$ws.Range("AQ77").Value = 255
$ws.Range("BO77").Value = 254
$ws.Range("BL78").Value = 1798
$ws.Range("BL79").Value = 1575
$ws.Range("BT79").Value = 1164
$ws.Range("AX80").Value = 2406
$ws.Range("BH80").Value = 2852
$ws.Range("BL80").Value = 1961
$ws.Range("AR81").Value = 8194
$ws.Range("AX81").Value = 4032
$ws.Range("BJ81").Value = 3696
$ws.Range("BL81").Value = 3267
$ws.Range("BT81").Value = 5009
$ws.Range("BU81").Value = 2266
